# The sheet holds one row (header) + 63 two-row "Primera"/"Segunda" pairs
# of weekly price data (rows 2-127, columns A:R). This commit adds one more
# week's worth of readings: every existing row's data shifts down by one
# pair (2 rows), a brand-new date (44462) is inserted for the now-current
# week (rows 12-13), and the oldest pair (previously rows 126-127) is
# preserved as two new rows appended at the bottom (128-129).
#
# Concretely: new_row[r] = old_row[r-2] for r = 14..129, and rows 12/13 only
# get a new Fecha (column D). Walk the rows bottom-up so a row's old
# contents are read before they get overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 18   # columns A..R

for ($r = 129; $r -ge 14; $r--) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $src = $ws.Cells.Item($r - 2, $c).Value2
        $ws.Cells.Item($r, $c).Value = $src
    }
}

# New current-week date for the top-most (now duplicated-forward) pair.
$ws.Cells.Item(12, 4).Value = 44462
$ws.Cells.Item(13, 4).Value = 44462

# The two brand-new rows need the date number format explicitly, since
# they don't inherit it the way the shifted-into rows do.
$ws.Cells.Item(128, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(129, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
